# "Generate Report for Handback"
#
# The handback-status report records a handful of timestamps tracking when
# files were handed off / handed back for localization. This run refreshes
# those timestamps to the newly generated values:
#
#   Overview!G2 (Latest HO Xliff Generate Date)     2016-08-16 02:57:30 -> 2016-08-16 02:58:19
#   zh-cn!H2    (Correspond Handoff Datetime)       2016-08-16 02:57:25 -> 2016-08-16 02:58:13
#   zh-cn!K2    (Correspond Handback DateTime)      2016-08-16 02:57:44 -> 2016-08-16 02:58:30
#   de-de!H2    (Correspond Handoff Datetime)       2016-08-16 02:57:30 -> 2016-08-16 02:58:19
#   de-de!K2    (Correspond Handback DateTime)      2016-08-16 02:57:51 -> 2016-08-16 02:58:37

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# These cells are formatted as dates (yyyy-mm-dd HH:mm:ss) but store the
# timestamp as literal text, so Value2 (which bypasses implicit date
# coercion) is used to keep the stored representation identical in kind.
$wsOverview.Range("G2").Value2 = "2016-08-16 02:58:19"

$wsZhCn.Range("H2").Value2 = "2016-08-16 02:58:13"
$wsZhCn.Range("K2").Value2 = "2016-08-16 02:58:30"

$wsDeDe.Range("H2").Value2 = "2016-08-16 02:58:19"
$wsDeDe.Range("K2").Value2 = "2016-08-16 02:58:37"
